$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine last used row from column A (Beteckning)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C ("Förändrad") holds a date serial value that gets bumped by
# one day (45180 -> 45181) for every data row (rows 2..lastRow) on this
# automatic refresh.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
